# Updates the cryptos list (Coin/Link/Price/Volume(1h)) on Sheet1
# to the latest scraped snapshot. Numeric-looking "Price" strings are
# written through a Text-formatted cell (then ClearFormats() restores the
# default, unstyled cell) so Excel stores them as literal text, exactly
# like the inline strings already used for values such as "1.00" or "207.71".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# Row 2: Bitcoin
$ws.Cells.Item(2, 4).Value = "27.080.05"
$ws.Cells.Item(2, 5).Value = "  +1.17%  "

# Row 3: Ethereum
$ws.Cells.Item(3, 4).Value = "1.567.73"
$ws.Cells.Item(3, 5).Value = "  +1.98%  "

# Row 4: TetherUSD
Set-TextValue 4 4 "1.00"
$ws.Cells.Item(4, 5).Value = "  +0.25%  "

# Row 5: BNB
Set-TextValue 5 4 "207.94"
$ws.Cells.Item(5, 5).Value = "  +1.28%  "

# Row 6: XRP
Set-TextValue 6 4 "0.490"
$ws.Cells.Item(6, 5).Value = "  +1.18%  "

# Row 7: USDC
$ws.Cells.Item(7, 5).Value = "  +0.36%  "

# Row 8: Solana
Set-TextValue 8 4 "22.01"
$ws.Cells.Item(8, 5).Value = "  +3.75%  "

# Row 9: Cardano
$ws.Cells.Item(9, 5).Value = "  +1.76%  "

# Row 10: Dogecoin
$ws.Cells.Item(10, 5).Value = "  +1.35%  "

# Row 11: TRON
$ws.Cells.Item(11, 5).Value = "  +0.54%  "

# Row 12: WrappedliquidstakedEther2.0
$ws.Cells.Item(12, 4).Value = "1.790.40"
$ws.Cells.Item(12, 5).Value = "  +1.96%  "

# Row 13: WrappedEther
$ws.Cells.Item(13, 4).Value = "1.566.95"
$ws.Cells.Item(13, 5).Value = "  +1.90%  "

# Row 14: Polkadot
$ws.Cells.Item(14, 5).Value = "  +2.59%  "

# Row 15: Polygon
Set-TextValue 15 4 "0.521"
$ws.Cells.Item(15, 5).Value = "  +2.67%  "

# Row 16: WrappedBTC
$ws.Cells.Item(16, 4).Value = "27.074.68"
$ws.Cells.Item(16, 5).Value = "  +1.15%  "

# Row 17: Litecoin
Set-TextValue 17 4 "62.00"
$ws.Cells.Item(17, 5).Value = "  +1.81%  "

# Row 18: BitcoinCash
Set-TextValue 18 4 "218.43"
$ws.Cells.Item(18, 5).Value = "  +2.59%  "

# Row 19: ShibaInu
$ws.Cells.Item(19, 4).Value = "0.0₃0698"
$ws.Cells.Item(19, 5).Value = "  +2.50%  "

# Row 20: Chainlink
$ws.Cells.Item(20, 5).Value = "  +1.46%  "

# Row 21: Dai
$ws.Cells.Item(21, 5).Value = "  +0.21%  "

# Row 22: Uniswap
Set-TextValue 22 4 "4.08"
$ws.Cells.Item(22, 5).Value = "  +1.80%  "

# Row 23: Avalanche
$ws.Cells.Item(23, 5).Value = "  +1.65%  "

# Row 24: Toncoin
$ws.Cells.Item(24, 5).Value = "  +1.39%  "

# Row 25: Monero
Set-TextValue 25 4 "154.35"
$ws.Cells.Item(25, 5).Value = "  +1.60%  "

# Row 26: Cosmos
$ws.Cells.Item(26, 5).Value = "  +0.81%  "

# Row 27: EthereumClassic
Set-TextValue 27 4 "14.97"
$ws.Cells.Item(27, 5).Value = "  +1.23%  "

# Row 28: BinanceUSD
$ws.Cells.Item(28, 5).Value = "  +0.34%  "

# Row 29: Stellar
$ws.Cells.Item(29, 5).Value = "  +1.66%  "

# Row 30: Hedera
Set-TextValue 30 4 "0.0471"
$ws.Cells.Item(30, 5).Value = "  +3.25%  "

# Row 31: PancakeSwap
$ws.Cells.Item(31, 5).Value = "  +0.53%  "

# Row 32: Filecoin
$ws.Cells.Item(32, 5).Value = "  +0.64%  "

# Row 33: Maker
$ws.Cells.Item(33, 4).Value = "1.447.46"
$ws.Cells.Item(33, 5).Value = "  +6.01%  "

# Row 34: InternetComputer(DFINITY)
Set-TextValue 34 4 "3.05"
$ws.Cells.Item(34, 5).Value = "  +4.42%  "

# Row 35: LidoDAOToken
$ws.Cells.Item(35, 5).Value = "  +4.47%  "

# Row 36: TrustWalletToken
Set-TextValue 36 4 "0.963"

# Row 37: HuobiToken
$ws.Cells.Item(37, 5).Value = "  +0.88%  "

# Row 38: VeChain
$ws.Cells.Item(38, 5).Value = "  +0.66%  "

# Row 39: ImmutableX
Set-TextValue 39 4 "0.523"
$ws.Cells.Item(39, 5).Value = "  +0.66%  "

# Row 40: ARBITRUM
Set-TextValue 40 4 "0.816"
$ws.Cells.Item(40, 5).Value = "  +1.71%  "

# Row 41: FraxShare
$ws.Cells.Item(41, 5).Value = "  +0.40%  "

# Row 42: PaxDollar
$ws.Cells.Item(42, 2).Value = "PaxDollar"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue 42 4 "1.01"
$ws.Cells.Item(42, 5).Value = "  +0.36%  "

# Row 43: WEMIXToken
$ws.Cells.Item(43, 2).Value = "WEMIXToken"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue 43 4 "0.992"
$ws.Cells.Item(43, 5).Value = "  -0.28%  "

# Row 44: MXToken
$ws.Cells.Item(44, 2).Value = "MXToken"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue 44 4 "2.28"
$ws.Cells.Item(44, 5).Value = "  +3.88%  "

# Row 45: Aave
Set-TextValue 45 4 "64.47"
$ws.Cells.Item(45, 5).Value = "  +2.58%  "

# Row 46: RenderToken
Set-TextValue 46 4 "1.76"
$ws.Cells.Item(46, 5).Value = "  +1.65%  "

# Row 47: RocketPoolETH
$ws.Cells.Item(47, 4).Value = "1.703.71"
$ws.Cells.Item(47, 5).Value = "  +1.97%  "

# Row 48: Quant
$ws.Cells.Item(48, 5).Value = "  +2.15%  "

# Row 49: Cronos
$ws.Cells.Item(49, 5).Value = "  +2.97%  "

# Row 50: Algorand
$ws.Cells.Item(50, 5).Value = "  +2.54%  "

# Row 51: USDD
Set-TextValue 51 4 "1.01"
$ws.Cells.Item(51, 5).Value = "  +0.62%  "

